$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("G2").Value = 2.05
$ws.Range("H2").Value = 3.4
$ws.Range("I2").Value = 3.6
$ws.Range("J2").Value = 2.75
$ws.Range("L2").Value = 4
$ws.Range("W2").Value = 7
$ws.Range("X2").Value = 9.5
$ws.Range("Y2").Value = 9
$ws.Range("Z2").Value = 19
$ws.Range("AA2").Value = 17
$ws.Range("AH2").Value = 9.5
$ws.Range("AI2").Value = 17
$ws.Range("AJ2").Value = 13
$ws.Range("AK2").Value = 41
$ws.Range("AL2").Value = 29
$ws.Range("AM2").Value = 41
$ws.Range("AN2").Value = 4
$ws.Range("AO2").Value = 12
$ws.Range("AR2").Value = 51
$ws.Range("AW2").Value = 5.5
$ws.Range("AX2").Value = 21
$ws.Range("AZ2").Value = 67
$ws.Range("BA2").Value = 101

# Row 3 updates
$ws.Range("Q3").Value = 2.1
$ws.Range("R3").Value = 1.7
